$wb = $excel.ActiveWorkbook

# Sheet "2025" (xl/worksheets/sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003837439598571561
$ws.Range("E2").Value = 0.371816737837252
$ws.Range("G2").Value = 0.2494892361375054
$ws.Range("I2").Value = 0.3688729365116042
$ws.Range("L2").Value = 0.597153
$ws.Range("M2").Value = 0.0822565
$ws.Range("N2").Value = 12.82009457445574
$ws.Range("O2").Value = 3.537862477780193

# Sheet "2030" (xl/worksheets/sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.009260248310130134
$ws.Range("B2").Value = 0.04893136040142861
$ws.Range("E2").Value = 0.221690802927771
$ws.Range("I2").Value = 0.4200935979461164
$ws.Range("L2").Value = 0.116304098642403
$ws.Range("M2").Value = 0.04737166666666669
$ws.Range("N2").Value = 5.019459357782987
$ws.Range("O2").Value = 2.350434537569127

# Sheet "2035" (xl/worksheets/sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08320532380150918
$ws.Range("B2").Value = 0.02828327743582346
$ws.Range("E2").Value = 0.1707960241217166
$ws.Range("I2").Value = 0.5131155002332095
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04895450000000007
$ws.Range("N2").Value = 8.420345971561265
$ws.Range("O2").Value = 3.507480209960086

# Sheet "2045" (xl/worksheets/sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("O2").Value = 4.778372589157264
